$d = $word.ActiveDocument

$d.Content.Find.Execute("Type of transfer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "What type of transfer is it?", 2)
